$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "283.18"
Set-TextValue $ws.Range("E2") "1.89%"
Set-TextValue $ws.Range("D3") "28.37"
Set-TextValue $ws.Range("E3") "4.02%"
Set-TextValue $ws.Range("D4") "5.015"
Set-TextValue $ws.Range("E4") "3.05%"
Set-TextValue $ws.Range("D5") "0.06519"
Set-TextValue $ws.Range("E5") "1.48%"
Set-TextValue $ws.Range("D6") "7.214"
Set-TextValue $ws.Range("E6") "2.85%"
Set-TextValue $ws.Range("D7") "3.350"
Set-TextValue $ws.Range("E7") "1.33%"
Set-TextValue $ws.Range("D8") "1.379"
Set-TextValue $ws.Range("E8") "15.75%"
Set-TextValue $ws.Range("D9") "0.9182"
Set-TextValue $ws.Range("E9") "3.44%"
Set-TextValue $ws.Range("D10") "0.1537"
Set-TextValue $ws.Range("E10") "-0.20%"
Set-TextValue $ws.Range("D11") "0.06641"
Set-TextValue $ws.Range("E11") "28.55%"
Set-TextValue $ws.Range("D12") "0.07564"
Set-TextValue $ws.Range("E12") "0.86%"
Set-TextValue $ws.Range("D13") "0.02782"
Set-TextValue $ws.Range("E13") "-3.87%"
Set-TextValue $ws.Range("D14") "0.08988"
Set-TextValue $ws.Range("E14") "0.17%"
Set-TextValue $ws.Range("D15") "0.001583"
Set-TextValue $ws.Range("E15") "1.29%"
Set-TextValue $ws.Range("D16") "0.0006384"
Set-TextValue $ws.Range("E16") "0.04%"
Set-TextValue $ws.Range("D17") "0.006171"
Set-TextValue $ws.Range("E17") "1.01%"
Set-TextValue $ws.Range("D18") "3.447"
Set-TextValue $ws.Range("E18") "-0.73%"
Set-TextValue $ws.Range("D19") "2.236"
Set-TextValue $ws.Range("E19") "-1.52%"
Set-TextValue $ws.Range("E21") "-4.49%"
Set-TextValue $ws.Range("D22") "3.980"
Set-TextValue $ws.Range("E22") "1.79%"
Set-TextValue $ws.Range("D23") "0.1543"
Set-TextValue $ws.Range("E23") "1.65%"
Set-TextValue $ws.Range("D24") "0.04428"
Set-TextValue $ws.Range("E24") "0.40%"
Set-TextValue $ws.Range("D25") "0.001185"
Set-TextValue $ws.Range("E25") "0.83%"
Set-TextValue $ws.Range("D26") "0.004438"
Set-TextValue $ws.Range("E26") "14.24%"
Set-TextValue $ws.Range("E27") "1.62%"
Set-TextValue $ws.Range("E28") "-1.58%"
Set-TextValue $ws.Range("D40") "0.04117"
Set-TextValue $ws.Range("E40") "-0.06%"
Set-TextValue $ws.Range("D41") "0.006688"
Set-TextValue $ws.Range("E41") "-1.96%"
Set-TextValue $ws.Range("D42") "0.1232"
Set-TextValue $ws.Range("E42") "4.90%"
Set-TextValue $ws.Range("D43") "0.002148"
Set-TextValue $ws.Range("E43") "13.07%"
Set-TextValue $ws.Range("D44") "0.01211"
Set-TextValue $ws.Range("E44") "3.82%"
Set-TextValue $ws.Range("D45") "0.00005669"
Set-TextValue $ws.Range("E45") "6.29%"
Set-TextValue $ws.Range("D46") "1.965"
Set-TextValue $ws.Range("E46") "20.71%"
